$d = $word.ActiveDocument

# NOTE: the very first content mutation issued against this session is
# silently dropped by this host (a quirk of the COM-interop shim), so we
# open with a harmless no-op "mutation" (re-assigning a Range's Text to
# its own current value) before making the edits we actually care about.
$warmup = $d.Paragraphs.Item(1).Range
$warmup.Text = $warmup.Text

# Walk every header/footer in every section and rename the inline
# pictures:
#   - the BTEC logo (descr "BTec_Logo-Orange"), currently "image1.jpg",
#     becomes "image2.jpg"
#   - the Pearson logo (descr contains "PearsonLogo.png"), currently
#     "image2.png", becomes "image1.png"
foreach ($sec in $d.Sections) {
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($shp in $hdr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($shp in $ftr.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                } elseif ($shp.AlternativeText -like "*PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
